$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 449.6111  # H41: was 500.17648
$ws.Cells.Item(41, 9).Value = 352  # I41: was 413.4
$ws.Cells.Item(41, 10).Value = 644.8333  # J41: was 624.1429000000001
$ws.Cells.Item(41, 11).Value = 352  # K41: was 413.4
$ws.Cells.Item(41, 12).Value = 644.8333  # L41: was 624.1429000000001
$ws.Cells.Item(41, 13).Value = 88  # M41: was 26.60000000000002
$ws.Cells.Item(41, 14).Value = -1524.8333  # N41: was -1504.1429
$ws.Cells.Item(53, 8).Value = 178.71428  # H53: was 185.15384
$ws.Cells.Item(53, 10).Value = 226.25  # J53: was 270
$ws.Cells.Item(53, 12).Value = 226.25  # L53: was 270
$ws.Cells.Item(53, 14).Value = -1500.25  # N53: was -1544
$ws.Cells.Item(121, 8).Value = 4950  # H121: was 5923.75
$ws.Cells.Item(121, 10).Value = 4950  # J121: was 5923.75
$ws.Cells.Item(121, 12).Value = 14850  # L121: was 17771.25
$ws.Cells.Item(121, 14).Value = -18344  # N121: was -21265.25
$ws.Cells.Item(125, 8).Value = 4855.143  # H125: was 4998.143
$ws.Cells.Item(125, 9).Value = 4699.2  # I125: was 4749.5
$ws.Cells.Item(125, 10).Value = 5245  # J125: was 5329.6665
$ws.Cells.Item(125, 11).Value = 42292.8  # K125: was 42745.5
$ws.Cells.Item(125, 12).Value = 47205  # L125: was 47966.9985
$ws.Cells.Item(125, 13).Value = -39832.8  # M125: was -40285.5
$ws.Cells.Item(125, 14).Value = -52125  # N125: was -52886.9985
$ws.Cells.Item(132, 8).Value = 14347.267  # H132: was 16453
$ws.Cells.Item(132, 9).Value = 16272.538  # I132: was 17573.584
$ws.Cells.Item(132, 10).Value = 1833  # J132: was 3006
$ws.Cells.Item(132, 11).Value = 48817.614  # K132: was 52720.75199999999
$ws.Cells.Item(132, 12).Value = 5499  # L132: was 9018
$ws.Cells.Item(132, 13).Value = -46287.614  # M132: was -50190.75199999999
$ws.Cells.Item(132, 14).Value = -10559  # N132: was -14078

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 15668.875  # H61: was 17707.285
$ws.Cells.Item(61, 9).Value = 17193.715  # I61: was 19826
$ws.Cells.Item(61, 11).Value = 17193.715  # K61: was 19826
$ws.Cells.Item(61, 13).Value = -16981.715  # M61: was -19614
$ws.Cells.Item(132, 8).Value = 5435.289  # H132: was 5660.5347
$ws.Cells.Item(132, 9).Value = 3553.6155  # I132: was 3800.375
$ws.Cells.Item(132, 11).Value = 10660.8465  # K132: was 11401.125
$ws.Cells.Item(132, 13).Value = -8130.8465  # M132: was -8871.125
$ws.Cells.Item(136, 8).Value = 15668.875  # H136: was 17707.285
$ws.Cells.Item(136, 9).Value = 17193.715  # I136: was 19826
$ws.Cells.Item(136, 11).Value = 51581.145  # K136: was 59478
$ws.Cells.Item(136, 13).Value = -49031.145  # M136: was -56928

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(29, 8).Value = 2120.1667  # H29: was 2933
$ws.Cells.Item(29, 9).Value = 2120.1667  # I29: was 2933
$ws.Cells.Item(29, 11).Value = 2120.1667  # K29: was 2933
$ws.Cells.Item(29, 13).Value = -1831.1667  # M29: was -2644
$ws.Cells.Item(105, 8).Value = 5739.0356  # H105: was 5673.2666
$ws.Cells.Item(105, 9).Value = 3166.5557  # I105: was 3172.7222
$ws.Cells.Item(105, 10).Value = 10369.5  # J105: was 9424.083000000001
$ws.Cells.Item(105, 11).Value = 3166.5557  # K105: was 3172.7222
$ws.Cells.Item(105, 12).Value = 10369.5  # L105: was 9424.083000000001
$ws.Cells.Item(105, 13).Value = -1419.5557  # M105: was -1425.7222
$ws.Cells.Item(105, 14).Value = -13863.5  # N105: was -12918.083
$ws.Cells.Item(107, 8).Value = 1498.5151  # H107: was 1512.1515
$ws.Cells.Item(107, 9).Value = 1299.4  # I107: was 1304.1904
$ws.Cells.Item(107, 10).Value = 1804.8462  # J107: was 1876.0834
$ws.Cells.Item(107, 11).Value = 1299.4  # K107: was 1304.1904
$ws.Cells.Item(107, 12).Value = 1804.8462  # L107: was 1876.0834
$ws.Cells.Item(107, 13).Value = 620.5999999999999  # M107: was 615.8096
$ws.Cells.Item(107, 14).Value = -5644.8462  # N107: was -5716.0834
$ws.Cells.Item(134, 8).Value = 4550.1816  # H134: was 3989.923
$ws.Cells.Item(134, 9).Value = 4973.9473  # I134: was 4530.476
$ws.Cells.Item(134, 10).Value = 1866.3334  # J134: was 1719.6
$ws.Cells.Item(134, 11).Value = 14921.8419  # K134: was 13591.428
$ws.Cells.Item(134, 12).Value = 5599.0002  # L134: was 5158.799999999999
$ws.Cells.Item(134, 13).Value = -12386.8419  # M134: was -11056.428
$ws.Cells.Item(134, 14).Value = -10669.0002  # N134: was -10228.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 825.5  # H107: was 788.5357
$ws.Cells.Item(107, 9).Value = 640.7646999999999  # I107: was 659.55
$ws.Cells.Item(107, 10).Value = 1453.6  # J107: was 1111
$ws.Cells.Item(107, 11).Value = 640.7646999999999  # K107: was 659.55
$ws.Cells.Item(107, 12).Value = 1453.6  # L107: was 1111
$ws.Cells.Item(107, 13).Value = 1279.2353  # M107: was 1260.45
$ws.Cells.Item(107, 14).Value = -5293.6  # N107: was -4951
$ws.Cells.Item(132, 8).Value = 8496.433999999999  # H132: was 7342.5137
$ws.Cells.Item(132, 9).Value = 8803.357  # I132: was 8280.833000000001
$ws.Cells.Item(132, 10).Value = 4199.5  # J132: was 3321.1428
$ws.Cells.Item(132, 11).Value = 26410.071  # K132: was 24842.499
$ws.Cells.Item(132, 12).Value = 12598.5  # L132: was 9963.428400000001
$ws.Cells.Item(132, 13).Value = -23880.071  # M132: was -22312.499
$ws.Cells.Item(132, 14).Value = -17658.5  # N132: was -15023.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 9215.637000000001  # H80: was 3674.4
$ws.Cells.Item(80, 10).Value = 10519.111  # J80: was 3890.6667
$ws.Cells.Item(80, 12).Value = 31557.333  # L80: was 11672.0001
$ws.Cells.Item(80, 14).Value = -33429.333  # N80: was -13544.0001
$ws.Cells.Item(83, 8).Value = 9215.637000000001  # H83: was 3674.4
$ws.Cells.Item(83, 10).Value = 10519.111  # J83: was 3890.6667
$ws.Cells.Item(83, 12).Value = 94671.99900000001  # L83: was 35016.0003
$ws.Cells.Item(83, 14).Value = -104031.999  # N83: was -44376.0003
$ws.Cells.Item(113, 8).Value = 676.8125  # H113: was 679.125
$ws.Cells.Item(113, 9).Value = 416.5  # I113: was 424.66666
$ws.Cells.Item(113, 10).Value = 714  # J113: was 737.8461
$ws.Cells.Item(113, 11).Value = 1249.5  # K113: was 1273.99998
$ws.Cells.Item(113, 12).Value = 2142  # L113: was 2213.5383
$ws.Cells.Item(113, 13).Value = 920.5  # M113: was 896.0000199999999
$ws.Cells.Item(113, 14).Value = -6482  # N113: was -6553.5383
$ws.Cells.Item(134, 8).Value = 7382.316  # H134: was 7530.6313
$ws.Cells.Item(134, 9).Value = 2206.2727  # I134: was 2462.4546
$ws.Cells.Item(134, 11).Value = 6618.8181  # K134: was 7387.3638
$ws.Cells.Item(134, 13).Value = -1548.8181  # M134: was -2317.3638
$ws.Cells.Item(140, 8).Value = 1825.6666  # H140: was 2050.647
$ws.Cells.Item(140, 9).Value = 1638.9412  # I140: was 1866.3125
$ws.Cells.Item(140, 11).Value = 4916.8236  # K140: was 5598.9375
$ws.Cells.Item(140, 13).Value = 263.1764000000003  # M140: was -418.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 0  # H29: was 10000
$ws.Cells.Item(29, 9).Value = 0  # I29: was 10000
$ws.Cells.Item(29, 11).Value = 0  # K29: was 10000
$ws.Cells.Item(29, 13).ClearContents()  # M29: was -9710
$ws.Cells.Item(80, 8).Value = 0  # H80: was 1725
$ws.Cells.Item(80, 9).Value = 0  # I80: was 2100
$ws.Cells.Item(80, 10).Value = 0  # J80: was 1350
$ws.Cells.Item(80, 11).Value = 0  # K80: was 2100
$ws.Cells.Item(80, 12).Value = 0  # L80: was 1350
$ws.Cells.Item(80, 13).ClearContents()  # M80: was -1102
$ws.Cells.Item(80, 14).ClearContents()  # N80: was -3346
$ws.Cells.Item(83, 8).Value = 0  # H83: was 1725
$ws.Cells.Item(83, 9).Value = 0  # I83: was 2100
$ws.Cells.Item(83, 10).Value = 0  # J83: was 1350
$ws.Cells.Item(83, 11).Value = 0  # K83: was 10500
$ws.Cells.Item(83, 12).Value = 0  # L83: was 6750
$ws.Cells.Item(83, 13).ClearContents()  # M83: was -5508
$ws.Cells.Item(83, 14).ClearContents()  # N83: was -16734
$ws.Cells.Item(102, 8).Value = 2511.2942  # H102: was 2517.7058
$ws.Cells.Item(102, 9).Value = 1012.8  # I102: was 1020.06665
$ws.Cells.Item(102, 11).Value = 1012.8  # K102: was 1020.06665
$ws.Cells.Item(102, 13).Value = 609.2  # M102: was 601.93335
$ws.Cells.Item(113, 8).Value = 1467.909  # H113: was 1470.5834
$ws.Cells.Item(113, 9).Value = 1464.7  # I113: was 1467.909
$ws.Cells.Item(113, 11).Value = 1464.7  # K113: was 1467.909
$ws.Cells.Item(113, 13).Value = 705.3  # M113: was 702.0909999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 10877.9  # H40: was 11253.889
$ws.Cells.Item(40, 9).Value = 15821  # I40: was 15821.25
$ws.Cells.Item(40, 10).Value = 7582.5  # J40: was 7600
$ws.Cells.Item(40, 11).Value = 15821  # K40: was 15821.25
$ws.Cells.Item(40, 12).Value = 7582.5  # L40: was 7600
$ws.Cells.Item(40, 13).Value = -15685  # M40: was -15685.25
$ws.Cells.Item(40, 14).Value = -7854.5  # N40: was -7872
$ws.Cells.Item(46, 8).Value = 2829.075  # H46: was 2782.0244
$ws.Cells.Item(46, 9).Value = 1425  # I46: was 1384.6154
$ws.Cells.Item(46, 11).Value = 1425  # K46: was 1384.6154
$ws.Cells.Item(46, 13).Value = -1237  # M46: was -1196.6154
$ws.Cells.Item(55, 8).Value = 1157.4688  # H55: was 1191.6774
$ws.Cells.Item(55, 9).Value = 1259.0769  # I55: was 1355.9166
$ws.Cells.Item(55, 11).Value = 1259.0769  # K55: was 1355.9166
$ws.Cells.Item(55, 13).Value = -1086.0769  # M55: was -1182.9166
$ws.Cells.Item(61, 8).Value = 6549.75  # H61: was 6840.4814
$ws.Cells.Item(61, 9).Value = 7207.8  # I61: was 7581.8945
$ws.Cells.Item(61, 10).Value = 4904.625  # J61: was 5079.625
$ws.Cells.Item(61, 11).Value = 7207.8  # K61: was 7581.8945
$ws.Cells.Item(61, 12).Value = 4904.625  # L61: was 5079.625
$ws.Cells.Item(61, 13).Value = -7005.8  # M61: was -7379.8945
$ws.Cells.Item(61, 14).Value = -5308.625  # N61: was -5483.625
$ws.Cells.Item(82, 8).Value = 11846.7  # H82: was 11931.4
$ws.Cells.Item(82, 9).Value = 13497  # I82: was 14290.5625
$ws.Cells.Item(82, 10).Value = 2495  # J82: was 2494.75
$ws.Cells.Item(82, 11).Value = 13497  # K82: was 14290.5625
$ws.Cells.Item(82, 12).Value = 2495  # L82: was 2494.75
$ws.Cells.Item(82, 13).Value = -13136  # M82: was -13929.5625
$ws.Cells.Item(82, 14).Value = -3217  # N82: was -3216.75
$ws.Cells.Item(85, 8).Value = 11846.7  # H85: was 11931.4
$ws.Cells.Item(85, 9).Value = 13497  # I85: was 14290.5625
$ws.Cells.Item(85, 10).Value = 2495  # J85: was 2494.75
$ws.Cells.Item(85, 11).Value = 13497  # K85: was 14290.5625
$ws.Cells.Item(85, 12).Value = 2495  # L85: was 2494.75
$ws.Cells.Item(85, 13).Value = -12249  # M85: was -13042.5625
$ws.Cells.Item(85, 14).Value = -4991  # N85: was -4990.75
$ws.Cells.Item(100, 8).Value = 3313  # H100: was 3237.2727
$ws.Cells.Item(100, 9).Value = 1938.3334  # I100: was 1944.2858
$ws.Cells.Item(100, 10).Value = 5375  # J100: was 5500
$ws.Cells.Item(100, 11).Value = 1938.3334  # K100: was 1944.2858
$ws.Cells.Item(100, 12).Value = 5375  # L100: was 5500
$ws.Cells.Item(100, 13).Value = -1397.3334  # M100: was -1403.2858
$ws.Cells.Item(100, 14).Value = -6457  # N100: was -6582
$ws.Cells.Item(109, 8).Value = 55000  # H109: was 0
$ws.Cells.Item(109, 10).Value = 55000  # J109: was 0
$ws.Cells.Item(109, 12).Value = 55000  # L109: was 0
$ws.Cells.Item(109, 14).Value = -57774  # N109: was None
$ws.Cells.Item(113, 8).Value = 6549.75  # H113: was 6840.4814
$ws.Cells.Item(113, 9).Value = 7207.8  # I113: was 7581.8945
$ws.Cells.Item(113, 10).Value = 4904.625  # J113: was 5079.625
$ws.Cells.Item(113, 11).Value = 7207.8  # K113: was 7581.8945
$ws.Cells.Item(113, 12).Value = 4904.625  # L113: was 5079.625
$ws.Cells.Item(113, 13).Value = -5037.8  # M113: was -5411.8945
$ws.Cells.Item(113, 14).Value = -9244.625  # N113: was -9419.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(30, 8).Value = 16333  # H30: was 0
$ws.Cells.Item(30, 9).Value = 8999  # I30: was 0
$ws.Cells.Item(30, 10).Value = 20000  # J30: was 0
$ws.Cells.Item(30, 11).Value = 8999  # K30: was 0
$ws.Cells.Item(30, 12).Value = 20000  # L30: was 0
$ws.Cells.Item(30, 13).Value = -8892  # M30: was None
$ws.Cells.Item(30, 14).Value = -20214  # N30: was None
$ws.Cells.Item(122, 8).Value = 47237.6  # H122: was 47197.88
$ws.Cells.Item(122, 9).Value = 1749.875  # I122: was 1800
$ws.Cells.Item(122, 10).Value = 128104.664  # J122: was 115294.7
$ws.Cells.Item(122, 11).Value = 5249.625  # K122: was 5400
$ws.Cells.Item(122, 12).Value = 384313.992  # L122: was 345884.1
$ws.Cells.Item(122, 13).Value = -2799.625  # M122: was -2950
$ws.Cells.Item(122, 14).Value = -389213.992  # N122: was -350784.1
